$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "293.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.13%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.22%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.152"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.51%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.530"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.89%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.616"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.400"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.11%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9176"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.97%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1637"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.80%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07831"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "21.93%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07753"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.72%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02945"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.30%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09005"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.19%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001588"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.69%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006566"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.79%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006189"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.18%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.484"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.10%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.243"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.31%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.12%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.08%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.151"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.10%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.08%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04536"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.72%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.80%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004238"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.43%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.45%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001689"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "4.40%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04415"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.73%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007043"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.34%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1275"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.16%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002208"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "11.56%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01325"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005840"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.66%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.726"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-12.42%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.58%"
